# Update the "ID" column (B) from hex-string message IDs to their
# plain decimal numeric equivalents (e.g. "0x10" -> 10), matching the
# "back up to date" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 10
    3  = 11
    4  = 12
    5  = 13
    6  = 14
    7  = 15
    8  = 16
    9  = 17
    10 = 18
    11 = 19
    12 = 50
    13 = 99
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row]
}

# Move the active selection (as recorded in the sheet view) to F15.
$ws.Range("F15").Select()
